$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows (14..104) down to (15..105)
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new data record
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100112031
$ws.Cells.Item(14, 7).Value = "Poroto verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 1400
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 1600
$ws.Cells.Item(14, 13).Value = 1546
$ws.Cells.Item(14, 14).Value = "$/kilo"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 1546
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"
